{"js": "const replacements = [\n  [\"2024-04-15 Monday\", \"2024-04-16 Tuesday\"],\n  [\"66\u00d770=\", \"69\u00d728=\"],\n  [\"82\u00d752=\", \"17\u00d718=\"],\n  [\"32\u00d749=\", \"22\u00d745=\"],\n  [\"26\u00d760=\", \"79\u00d744=\"],\n  [\"54\u00d725=\", \"21\u00d763=\"],\n  [\"44\u00d758=\", \"21\u00d764=\"],\n  [\"16\u00d787=\", \"98\u00d730=\"],\n  [\"81\u00d738=\", \"23\u00d789=\"],\n  [\"91\u00d786=\", \"69\u00d738=\"],\n  [\"96\u00d719=\", \"50\u00d737=\"],\n  [\"38\u00d787=\", \"44\u00d736=\"],\n  [\"15\u00d792=\", \"32\u00d735=\"],\n  [\"77\u00d764=\", \"72\u00d796=\"],\n  [\"13\u00d763=\", \"92\u00d787=\"],\n  [\"41\u00d759=\", \"81\u00d726=\"],\n  [\"63\u00d740=\", \"49\u00d758=\"],\n  [\"20\u00d740=\", \"19\u00d767=\"],\n  [\"80\u00d778=\", \"27\u00d779=\"],\n  [\"22\u00d789=\", \"34\u00d799=\"],\n  [\"17\u00d713=\", \"78\u00d740=\"],\n  [\"53\u00d738=\", \"90\u00d728=\"],\n  [\"69\u00d737=\", \"34\u00d740=\"],\n  [\"78\u00d745=\", \"69\u00d767=\"],\n  [\"12\u00d740=\", \"33\u00d730=\"],\n  [\"40\u00d739=\", \"49\u00d765=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-15 Monday\", \"2024-04-16 Tuesday\"),\n    @(\"66\u00d770=\", \"69\u00d728=\"),\n    @(\"82\u00d752=\", \"17\u00d718=\"),\n    @(\"32\u00d749=\", \"22\u00d745=\"),\n    @(\"26\u00d760=\", \"79\u00d744=\"),\n    @(\"54\u00d725=\", \"21\u00d763=\"),\n    @(\"44\u00d758=\", \"21\u00d764=\"),\n    @(\"16\u00d787=\", \"98\u00d730=\"),\n    @(\"81\u00d738=\", \"23\u00d789=\"),\n    @(\"91\u00d786=\", \"69\u00d738=\"),\n    @(\"96\u00d719=\", \"50\u00d737=\"),\n    @(\"38\u00d787=\", \"44\u00d736=\"),\n    @(\"15\u00d792=\", \"32\u00d735=\"),\n    @(\"77\u00d764=\", \"72\u00d796=\"),\n    @(\"13\u00d763=\", \"92\u00d787=\"),\n    @(\"41\u00d759=\", \"81\u00d726=\"),\n    @(\"63\u00d740=\", \"49\u00d758=\"),\n    @(\"20\u00d740=\", \"19\u00d767=\"),\n    @(\"80\u00d778=\", \"27\u00d779=\"),\n    @(\"22\u00d789=\", \"34\u00d799=\"),\n    @(\"17\u00d713=\", \"78\u00d740=\"),\n    @(\"53\u00d738=\", \"90\u00d728=\"),\n    @(\"69\u00d737=\", \"34\u00d740=\"),\n    @(\"78\u00d745=\", \"69\u00d767=\"),\n    @(\"12\u00d740=\", \"33\u00d730=\"),\n    @(\"40\u00d739=\", \"49\u00d765=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $wdReplaceOne = 1\n    $wdFindContinue = 1\n    $found = $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceOne)\n    if (-not $found) {\n        throw \"Could not find text: $($pair[0])\"\n    }\n}"}
